$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Magnesium chloride unit price (row 8): update loading value and
# replace the formula-driven lower/upper bounds with their frozen values.
$ws.Range("E8").Value = 0.38
$ws.Range("G8").Value = 0.349
$ws.Range("I8").Value = 0.411

# Zinc sulfate unit price (row 9): same treatment.
$ws.Range("E9").Value = 0.795
$ws.Range("G9").Value = 0.657
$ws.Range("I9").Value = 0.931

# Mirror the author's final on-screen selection (rows 8:9 highlighted,
# active cell A8) left over from editing those two rows.
$ws.Range("A8:XFD9").Select()
